$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1616.9783
$ws.Range("I40").Value = 1539.4667
$ws.Range("J40").Value = 1958.9412
$ws.Range("K40").Value = 1539.4667
$ws.Range("L40").Value = 1958.9412
$ws.Range("M40").Value = -1364.4667
$ws.Range("N40").Value = -2308.9412

$ws.Range("H64").Value = 3711.2432
$ws.Range("I64").Value = 3772.52
$ws.Range("J64").Value = 3583.5833
$ws.Range("K64").Value = 3772.52
$ws.Range("L64").Value = 3583.5833
$ws.Range("M64").Value = -3524.52
$ws.Range("N64").Value = -4079.5833

$ws.Range("H67").Value = 3711.2432
$ws.Range("I67").Value = 3772.52
$ws.Range("J67").Value = 3583.5833
$ws.Range("K67").Value = 3772.52
$ws.Range("L67").Value = 3583.5833
$ws.Range("M67").Value = -2914.52
$ws.Range("N67").Value = -5299.5833

$ws.Range("H113").Value = 3050
$ws.Range("I113").Value = 3050
$ws.Range("K113").Value = 3050
$ws.Range("M113").Value = 204

$ws.Range("H125").Value = 4463.7896
$ws.Range("I125").Value = 1372
$ws.Range("J125").Value = 6712.364
$ws.Range("K125").Value = 12348
$ws.Range("L125").Value = 60411.276
$ws.Range("M125").Value = -9888
$ws.Range("N125").Value = -65331.276

$ws.Range("H134").Value = 26240
$ws.Range("J134").Value = 26240
$ws.Range("L134").Value = 26240
$ws.Range("N134").Value = -36380

$ws.Range("H138").Value = 5030.5
$ws.Range("I138").Value = 923.76666
$ws.Range("J138").Value = 11875.056
$ws.Range("K138").Value = 2771.29998
$ws.Range("L138").Value = 35625.16800000001
$ws.Range("M138").Value = 2368.70002
$ws.Range("N138").Value = -45905.16800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4860.463
$ws.Range("I32").Value = 3688.3333
$ws.Range("J32").Value = 8962.916999999999
$ws.Range("K32").Value = 3688.3333
$ws.Range("L32").Value = 8962.916999999999
$ws.Range("M32").Value = -3401.3333
$ws.Range("N32").Value = -9536.916999999999

$ws.Range("H61").Value = 4582
$ws.Range("I61").Value = 5237.107
$ws.Range("J61").Value = 913.4
$ws.Range("K61").Value = 5237.107
$ws.Range("L61").Value = 913.4
$ws.Range("M61").Value = -5025.107
$ws.Range("N61").Value = -1337.4

$ws.Range("H74").Value = 1475.6025
$ws.Range("I74").Value = 1458.9855
$ws.Range("J74").Value = 1603
$ws.Range("K74").Value = 1458.9855
$ws.Range("L74").Value = 1603
$ws.Range("M74").Value = -584.9855
$ws.Range("N74").Value = -3351

$ws.Range("H77").Value = 1475.6025
$ws.Range("I77").Value = 1458.9855
$ws.Range("J77").Value = 1603
$ws.Range("K77").Value = 7294.9275
$ws.Range("L77").Value = 8015
$ws.Range("M77").Value = -2926.9275
$ws.Range("N77").Value = -16751

$ws.Range("H132").Value = 3177.1191
$ws.Range("I132").Value = 1679.6
$ws.Range("J132").Value = 5379.353
$ws.Range("K132").Value = 5038.799999999999
$ws.Range("L132").Value = 16138.059
$ws.Range("M132").Value = -2508.799999999999
$ws.Range("N132").Value = -21198.059

$ws.Range("H136").Value = 4582
$ws.Range("I136").Value = 5237.107
$ws.Range("J136").Value = 913.4
$ws.Range("K136").Value = 15711.321
$ws.Range("L136").Value = 2740.2
$ws.Range("M136").Value = -13161.321
$ws.Range("N136").Value = -7840.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1465.5385
$ws.Range("I94").Value = 927.44446
$ws.Range("J94").Value = 2676.25
$ws.Range("K94").Value = 927.44446
$ws.Range("L94").Value = 2676.25
$ws.Range("M94").Value = -476.44446
$ws.Range("N94").Value = -3578.25

$ws.Range("I99").Value = 111112330
$ws.Range("K99").Value = 111112330
$ws.Range("M99").Value = -111110832

$ws.Range("H137").Value = 42709
$ws.Range("I137").Value = 42709
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 42709
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = ""
$ws.Range("M137").Value = -37609

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6994034.5
$ws.Range("I16").Value = 10989926
$ws.Range("J16").Value = 1225
$ws.Range("K16").Value = 10989926
$ws.Range("L16").Value = 1225
$ws.Range("M16").Value = -10989639
$ws.Range("N16").Value = -1799

$ws.Range("H31").Value = 4299.7617
$ws.Range("I31").Value = 920.0645
$ws.Range("J31").Value = 13824.363
$ws.Range("K31").Value = 920.0645
$ws.Range("L31").Value = 13824.363
$ws.Range("M31").Value = -625.0645
$ws.Range("N31").Value = -14414.363

$ws.Range("H34").Value = 4299.7617
$ws.Range("I34").Value = 920.0645
$ws.Range("J34").Value = 13824.363
$ws.Range("K34").Value = 920.0645
$ws.Range("L34").Value = 13824.363
$ws.Range("M34").Value = -718.0645
$ws.Range("N34").Value = -14228.363

$ws.Range("H113").Value = 6994034.5
$ws.Range("I113").Value = 10989926
$ws.Range("J113").Value = 1225
$ws.Range("K113").Value = 10989926
$ws.Range("L113").Value = 1225
$ws.Range("M113").Value = -10987756
$ws.Range("N113").Value = -5565

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 5882514
$ws.Range("J12").Value = 170
$ws.Range("L12").Value = 510
$ws.Range("N12").Value = -856

$ws.Range("H122").Value = 802.64703
$ws.Range("I122").Value = 659.36365
$ws.Range("J122").Value = 1065.3334
$ws.Range("K122").Value = 5934.27285
$ws.Range("L122").Value = 9588.000599999999
$ws.Range("M122").Value = -3484.27285
$ws.Range("N122").Value = -14488.0006

$ws.Range("H132").Value = 3427.5942
$ws.Range("I132").Value = 3000.5
$ws.Range("J132").Value = 3483.6064
$ws.Range("K132").Value = 27004.5
$ws.Range("L132").Value = 31352.4576
$ws.Range("M132").Value = -24474.5
$ws.Range("N132").Value = -36412.4576

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 6033.0415
$ws.Range("I126").Value = 8849.429
$ws.Range("J126").Value = 2090.1
$ws.Range("K126").Value = 26548.287
$ws.Range("L126").Value = 6270.299999999999
$ws.Range("M126").Value = -24078.287
$ws.Range("N126").Value = -11210.3

$ws.Range("H132").Value = 3256.75
$ws.Range("I132").Value = 3322.3
$ws.Range("J132").Value = 3191.2
$ws.Range("K132").Value = 9966.900000000001
$ws.Range("L132").Value = 9573.599999999999
$ws.Range("M132").Value = -7436.900000000001
$ws.Range("N132").Value = -14633.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 6000
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 6000
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 6000
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = -6224

$ws.Range("H7").Value = 73850.64
$ws.Range("I7").Value = 127050.5
$ws.Range("J7").Value = 2917.5
$ws.Range("K7").Value = 127050.5
$ws.Range("L7").Value = 2917.5
$ws.Range("M7").Value = -126938.5
$ws.Range("N7").Value = -3141.5

$ws.Range("H46").Value = 15874000
$ws.Range("I46").Value = 23810254
$ws.Range("J46").Value = 1492.8572
$ws.Range("K46").Value = 23810254
$ws.Range("L46").Value = 1492.8572
$ws.Range("M46").Value = -23810066
$ws.Range("N46").Value = -1868.8572

$ws.Range("H61").Value = 3191.0833
$ws.Range("I61").Value = 2899
$ws.Range("J61").Value = 3600
$ws.Range("K61").Value = 2899
$ws.Range("L61").Value = 3600
$ws.Range("M61").Value = -2697
$ws.Range("N61").Value = -4004

$ws.Range("H93").Value = 55578428
$ws.Range("I93").Value = 33975
$ws.Range("J93").Value = 166667330
$ws.Range("K93").Value = 33975
$ws.Range("L93").Value = 166667330
$ws.Range("M93").Value = -32727
$ws.Range("N93").Value = -166669826

$ws.Range("H113").Value = 3191.0833
$ws.Range("I113").Value = 2899
$ws.Range("J113").Value = 3600
$ws.Range("K113").Value = 2899
$ws.Range("L113").Value = 3600
$ws.Range("M113").Value = -729
$ws.Range("N113").Value = -7940

$ws.Range("H126").Value = 73850.64
$ws.Range("I126").Value = 127050.5
$ws.Range("J126").Value = 2917.5
$ws.Range("K126").Value = 381151.5
$ws.Range("L126").Value = 8752.5
$ws.Range("M126").Value = -378681.5
$ws.Range("N126").Value = -13692.5

$ws.Range("H132").Value = 17340652
$ws.Range("I132").Value = 20642634
$ws.Range("K132").Value = 61927902
$ws.Range("M132").Value = -61925372

$ws.Range("H136").Value = 4856.4893
$ws.Range("I136").Value = 4354.9116
$ws.Range("J136").Value = 6168.3076
$ws.Range("K136").Value = 13064.7348
$ws.Range("L136").Value = 18504.9228
$ws.Range("M136").Value = -10514.7348
$ws.Range("N136").Value = -23604.9228

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1959.625
$ws.Range("I122").Value = 1816.5
$ws.Range("J122").Value = 2045.5
$ws.Range("K122").Value = 5449.5
$ws.Range("L122").Value = 6136.5
$ws.Range("M122").Value = -2999.5
$ws.Range("N122").Value = -11036.5

$ws.Range("H132").Value = 1608.2858
$ws.Range("I132").Value = 899.05554
$ws.Range("J132").Value = 2884.9
$ws.Range("K132").Value = 2697.16662
$ws.Range("L132").Value = 8654.700000000001
$ws.Range("M132").Value = -167.16662
$ws.Range("N132").Value = -13714.7
